# Loca_Keys_Strings_DE_EN.xlsx - add new localized strings for the
# "Get Free Chips" modal (global_get-free-chips-modal_*) below the
# existing KeyValuePairs table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reference cells carrying the two "data row" styles already used in
#     the sheet, so new cells reuse the existing style indices instead of
#     Excel fabricating brand-new ones:
#       style "key"       -> fill highlight, no wrap   (column A, e.g. A2)
#       style "wrapvalue" -> fill highlight, wrap text (e.g. B14)
$keyFmtSource   = $ws.Range("A2")
$wrapFmtSource  = $ws.Range("B14")

function Set-KeyStyle($cell) {
    $keyFmtSource.Copy()
    $cell.PasteSpecial(-4122) # xlPasteFormats
}

function Set-WrapStyle($cell) {
    $wrapFmtSource.Copy()
    $cell.PasteSpecial(-4122) # xlPasteFormats
}

# --- Row 72: global_get-free-chips-modal_header ----------------------
Set-KeyStyle  $ws.Range("A72")
Set-WrapStyle $ws.Range("B72")
Set-KeyStyle  $ws.Range("C72")

$ws.Range("A72").Value = "global_get-free-chips-modal_header"
$ws.Range("B72").Value = "Refuel  "
$ws.Range("C72").Value = "Auftanken"

# --- Row 73: global_get-free-chips-modal_content ----------------------
Set-KeyStyle  $ws.Range("A73")
Set-WrapStyle $ws.Range("B73")
Set-WrapStyle $ws.Range("C73")

$ws.Range("A73").Value = "global_get-free-chips-modal_content"
$ws.Range("B73").Value = "Oh noes, it seems like you're running out of chips! But don't worry, here's a fresh batch of chips for you so you can continue playing!"
$ws.Range("C73").Value = "Oh nein, es scheint, als würden Ihnen die Chips ausgehen! Aber keine Sorge, hier ist eine neue Charge Chips für Sie, damit Sie weiterspielen können!"
$ws.Rows(73).RowHeight = 75

# --- Row 74: global_get-free-chips-modal_btn-txt ----------------------
Set-KeyStyle $ws.Range("A74")
Set-KeyStyle $ws.Range("B74")
Set-KeyStyle $ws.Range("C74")

$ws.Range("A74").Value = "global_get-free-chips-modal_btn-txt"
$ws.Range("B74").Value = "Get Your Free Chips"
$ws.Range("C74").Value = "Gratis Chips Holen"

$excel.CutCopyMode = $false

# --- Minor row-height relayout that accompanied this edit (wrapped rows
#     re-measured slightly taller) --------------------------------------
$ws.Rows(14).RowHeight = 90
$ws.Rows(20).RowHeight = 45
$ws.Rows(47).RowHeight = 75
$ws.Rows(68).RowHeight = 30
$ws.Rows(71).RowHeight = 60

# --- Default workbook style is now labelled "Normal" instead of
#     "Standard" (locale of the authoring Excel changed) ---------------
try {
    $wb.Styles("Standard").Name = "Normal"
} catch {
}

# --- Leave the cursor where the author left it ------------------------
$ws.Activate() | Out-Null
$ws.Range("C70").Select() | Out-Null
